# Update res_bus vm_pu results for the Case_4_99 (380 kV) run.
# Slack bus voltage setpoint changed from 1.05 to 1.02 pu, and all
# downstream bus voltage magnitudes were recomputed accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nRows = 24

$bf = New-Object "object[,]" $nRows,5
$bf[0,0] = 1.02
$bf[0,1] = 1.104787407992471
$bf[0,2] = 1.097397575286285
$bf[0,3] = 1.115642162878807
$bf[0,4] = 1.116247020282871
$bf[1,0] = 1.02
$bf[1,1] = 1.106788846496812
$bf[1,2] = 1.098977415667224
$bf[1,3] = 1.117556641507543
$bf[1,4] = 1.118053299041216
$bf[2,0] = 1.02
$bf[2,1] = 1.10807913132044
$bf[2,2] = 1.099994988149107
$bf[2,3] = 1.118790961459637
$bf[2,4] = 1.11921759704334
$bf[3,0] = 1.02
$bf[3,1] = 1.108620447811198
$bf[3,2] = 1.100421672574504
$bf[3,3] = 1.119308821231391
$bf[3,4] = 1.119706016018488
$bf[4,0] = 1.02
$bf[4,1] = 1.108711272320189
$bf[4,2] = 1.100493250652425
$bf[4,3] = 1.119395711358725
$bf[4,4] = 1.119787962659183
$bf[5,0] = 1.02
$bf[5,1] = 1.108086368781944
$bf[5,2] = 1.100000693834803
$bf[5,3] = 1.118797885215163
$bf[5,4] = 1.119224127425295
$bf[6,0] = 1.02
$bf[6,1] = 1.105464808327935
$bf[6,2] = 1.097932472752405
$bf[6,3] = 1.116290110649544
$bf[6,4] = 1.116858402549194
$bf[7,0] = 1.02
$bf[7,1] = 1.100807545318151
$bf[7,2] = 1.094251191842799
$bf[7,3] = 1.111835765062054
$bf[7,4] = 1.112654354184124
$bf[8,0] = 1.02
$bf[8,1] = 1.097675736882937
$bf[8,2] = 1.091770997874666
$bf[8,3] = 1.108840974597988
$bf[8,4] = 1.109826503669876
$bf[9,0] = 1.02
$bf[9,1] = 1.096312875548706
$bf[9,2] = 1.090690595640738
$bf[9,3] = 1.107537882543787
$bf[9,4] = 1.108595734953088
$bf[10,0] = 1.02
$bf[10,1] = 1.095805600395432
$bf[10,2] = 1.090288290796894
$bf[10,3] = 1.107052876681035
$bf[10,4] = 1.108137600403055
$bf[11,0] = 1.02
$bf[11,1] = 1.095914460554615
$bf[11,2] = 1.090374631990936
$bf[11,3] = 1.107156956873377
$bf[11,4] = 1.108235916255462
$bf[12,0] = 1.02
$bf[12,1] = 1.096270965581639
$bf[12,2] = 1.090657361444898
$bf[12,3] = 1.107497811956263
$bf[12,4] = 1.108557885400277
$bf[13,0] = 1.02
$bf[13,1] = 1.096490480319008
$bf[13,2] = 1.090831427759135
$bf[13,3] = 1.107707693368042
$bf[13,4] = 1.108756131565604
$bf[14,0] = 1.02
$bf[14,1] = 1.097766040278757
$bf[14,2] = 1.091842562383561
$bf[14,3] = 1.108927320758688
$bf[14,4] = 1.109908050890358
$bf[15,0] = 1.02
$bf[15,1] = 1.098564331710874
$bf[15,2] = 1.092475073195161
$bf[15,3] = 1.10969064645526
$bf[15,4] = 1.110628916205334
$bf[16,0] = 1.02
$bf[16,1] = 1.099029310738463
$bf[16,2] = 1.092843384247926
$bf[16,3] = 1.110135272820212
$bf[16,4] = 1.111048779631103
$bf[17,0] = 1.02
$bf[17,1] = 1.099187747105467
$bf[17,2] = 1.092968864067027
$bf[17,3] = 1.110286776554749
$bf[17,4] = 1.111191840385675
$bf[18,0] = 1.02
$bf[18,1] = 1.098478750139755
$bf[18,2] = 1.092407275263558
$bf[18,3] = 1.109608811989212
$bf[18,4] = 1.110551636964534
$bf[19,0] = 1.02
$bf[19,1] = 1.096166012911887
$bf[19,2] = 1.090574132342888
$bf[19,3] = 1.107397465900635
$bf[19,4] = 1.108463100570751
$bf[20,0] = 1.02
$bf[20,1] = 1.094705825431505
$bf[20,2] = 1.089415791425411
$bf[20,3] = 1.106001423475371
$bf[20,4] = 1.107144315975303
$bf[21,0] = 1.02
$bf[21,1] = 1.095480485024542
$bf[21,2] = 1.090030405192137
$bf[21,3] = 1.106742040216759
$bf[21,4] = 1.107843972307103
$bf[22,0] = 1.02
$bf[22,1] = 1.098517422778041
$bf[22,2] = 1.092437912151099
$bf[22,3] = 1.109645791341539
$bf[22,4] = 1.110586557995948
$bf[23,0] = 1.02
$bf[23,1] = 1.102016194273051
$bf[23,2] = 1.095207381786813
$bf[23,3] = 1.112991657105387
$bf[23,4] = 1.113745531341639
$ws.Range("B2:F25").Value = $bf

$inArr = New-Object "object[,]" $nRows,6
$inArr[0,0] = 1.065597497592041
$inArr[0,1] = 1.109549191486713
$inArr[0,2] = 1.100027417661076
$inArr[0,3] = 1.11822670241576
$inArr[0,4] = 1.118830085185469
$inArr[0,5] = 1.111124878605696
$inArr[1,0] = 1.066221899488423
$inArr[1,1] = 1.111216232618573
$inArr[1,2] = 1.101426860832831
$inArr[1,3] = 1.119963335752586
$inArr[1,4] = 1.120458871351878
$inArr[1,5] = 1.112794287127175
$inArr[2,0] = 1.066622195055447
$inArr[2,1] = 1.112289834320856
$inArr[2,2] = 1.102327206304959
$inArr[2,3] = 1.121082136000439
$inArr[2,4] = 1.121507857048806
$inArr[2,5] = 1.113869413467019
$inArr[3,0] = 1.066789594684114
$inArr[3,1] = 1.112739980535885
$inArr[3,2] = 1.102704488330175
$inArr[3,3] = 1.121551324249647
$inArr[3,4] = 1.121947686833039
$inArr[3,5] = 1.114320198941306
$inArr[4,0] = 1.066817650212779
$inArr[4,1] = 1.11281549255571
$inArr[4,2] = 1.102767764453625
$inArr[4,3] = 1.121630035901407
$inArr[4,4] = 1.122021468540707
$inArr[4,5] = 1.11439581819686
$inArr[5,0] = 1.066624435320787
$inArr[5,1] = 1.112295853866749
$inArr[5,2] = 1.102332252346834
$inArr[5,3] = 1.121088409826973
$inArr[5,4] = 1.121513738620814
$inArr[5,5] = 1.113875441561358
$inArr[6,0] = 1.065809295992791
$inArr[6,1] = 1.110113641826887
$inArr[6,2] = 1.100501451572599
$inArr[6,3] = 1.118814637284393
$inArr[6,4] = 1.119381578133509
$inArr[6,5] = 1.111690130530096
$inArr[7,0] = 1.064343893375747
$inArr[7,1] = 1.106228384522181
$inArr[7,2] = 1.097234764976508
$inArr[7,3] = 1.114769304839532
$inArr[7,4] = 1.115585606617925
$inArr[7,5] = 1.107799355714372
$inArr[8,0] = 1.06334684774296
$inArr[8,1] = 1.103610002076866
$inArr[8,2] = 1.095028492784673
$inArr[8,3] = 1.112045030914976
$inArr[8,4] = 1.113027529523553
$inArr[8,5] = 1.105177254865652
$inArr[9,0] = 1.062910211641236
$inArr[9,1] = 1.102469214527894
$inArr[9,2] = 1.094066128427036
$inArr[9,3] = 1.110858581097167
$inArr[9,4] = 1.11191305130148
$inArr[9,5] = 1.104034847267511
$inArr[10,0] = 1.062747276362938
$inArr[10,1] = 1.102044395246649
$inArr[10,2] = 1.093707583335472
$inArr[10,3] = 1.110416828774142
$inArr[10,4] = 1.111498034986071
$inArr[10,5] = 1.103609424694158
$inArr[11,0] = 1.062782260625434
$inArr[11,1] = 1.10213556972974
$inArr[11,2] = 1.093784541731094
$inArr[11,3] = 1.11051163419518
$inArr[11,4] = 1.111587105303553
$inArr[11,5] = 1.103700728655473
$inArr[12,0] = 1.062896758711893
$inArr[12,1] = 1.102434121021122
$inArr[12,2] = 1.094036513148404
$inArr[12,3] = 1.110822087332484
$inArr[12,4] = 1.111878767482936
$inArr[12,5] = 1.103999703923934
$inArr[13,0] = 1.062967205133807
$inArr[13,1] = 1.102617924409445
$inArr[13,2] = 1.094191617223645
$inArr[13,3] = 1.111013227530408
$inArr[13,4] = 1.112058330333763
$inArr[13,5] = 1.10418376833415
$inArr[14,0] = 1.063375721426298
$inArr[14,1] = 1.103685562411366
$inArr[14,2] = 1.095092211555696
$inArr[14,3] = 1.112123625652376
$inArr[14,4] = 1.113101348020519
$inArr[14,5] = 1.105252922504494
$inArr[15,0] = 1.063630650484327
$inArr[15,1] = 1.104353367833394
$inArr[15,2] = 1.095655231029436
$inArr[15,3] = 1.11281830373361
$inArr[15,4] = 1.113753763105827
$inArr[15,5] = 1.105921676286819
$inArr[16,0] = 1.063778873466046
$inArr[16,1] = 1.104742212834356
$inArr[16,2] = 1.09598295368471
$inArr[16,3] = 1.113222841414668
$inArr[16,4] = 1.114133650098388
$inArr[16,5] = 1.106311073492284
$inArr[17,0] = 1.063829333824485
$inArr[17,1] = 1.104874685309112
$inArr[17,2] = 1.096094584498142
$inArr[17,3] = 1.113360667809059
$inArr[17,4] = 1.114263071303256
$inArr[17,5] = 1.106443734093149
$inArr[18,0] = 1.06360334799973
$inArr[18,1] = 1.10428178860134
$inArr[18,2] = 1.095594894569697
$inArr[18,3] = 1.112743839429008
$inArr[18,4] = 1.11368383314794
$inArr[18,5] = 1.105849995404045
$inArr[19,0] = 1.062863062655396
$inArr[19,1] = 1.102346235136662
$inArr[19,2] = 1.093962343852122
$inArr[19,3] = 1.110730695883348
$inArr[19,4] = 1.111792909420268
$inArr[19,5] = 1.103911693231437
$inArr[20,0] = 1.06239327537086
$inArr[20,1] = 1.101123013819363
$inArr[20,2] = 1.092929633352667
$inArr[20,3] = 1.109458852104066
$inArr[20,4] = 1.110597925597439
$inArr[20,5] = 1.102686734799561
$inArr[21,0] = 1.062642733884047
$inArr[21,1] = 1.1017720691787
$inArr[21,2] = 1.09347769409406
$inArr[21,3] = 1.110133668011982
$inArr[21,4] = 1.111231994625499
$inArr[21,5] = 1.103336711891957
$inArr[22,0] = 1.063615686277929
$inArr[22,1] = 1.104314134258843
$inArr[22,2] = 1.095622160110085
$inArr[22,3] = 1.112777488670018
$inArr[22,4] = 1.113715433511039
$inArr[22,5] = 1.105882386996092
$inArr[23,0] = 1.064726235370282
$inArr[23,1] = 1.107237687356702
$inArr[23,2] = 1.098084212990114
$inArr[23,3] = 1.115819843544714
$inArr[23,4] = 1.116571692868012
$inArr[23,5] = 1.108810091874687
$ws.Range("I2:N25").Value = $inArr

